# Updates cryptos list figures (price & 1h volume change) per the Tue Nov 14 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.614.69"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.059.94"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'243.50"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'0.669"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'54.65"
$ws.Range("E8").Value = "  -6.47%  "
$ws.Range("D9").Value = "'58.68"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "'0.941"
$ws.Range("E13").Value = "  +7.19%  "
$ws.Range("D14").Value = "'14.81"
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("D15").Value = "2.359.46"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'5.45"
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("D17").Value = "2.055.89"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "36.523.75"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'16.83"
$ws.Range("E19").Value = "  -7.42%  "
$ws.Range("D20").Value = "'72.05"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "'238.38"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'5.27"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").Value = "'2.14"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").Value = "'9.33"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").Value = "'20.15"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("E31").Value = "  +8.59%  "
$ws.Range("D32").Value = "'5.07"
$ws.Range("E32").Value = "  -6.76%  "
$ws.Range("D33").Value = "'4.50"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("D34").Value = "'0.0599"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "'2.21"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("E38").Value = "  -5.00%  "
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("D41").Value = "'0.0217"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").Value = "'2.86"
$ws.Range("E42").Value = "  -7.87%  "
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").Value = "'94.56"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0912"
$ws.Range("E45").Value = "  -5.41%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.419.72"
$ws.Range("E46").Value = "  +9.81%  "
$ws.Range("D47").Value = "'16.02"
$ws.Range("E47").Value = "  -5.48%  "
$ws.Range("D48").Value = "'7.49"
$ws.Range("E48").Value = "  +11.86%  "
$ws.Range("D49").Value = "'2.87"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "2.246.68"
$ws.Range("E51").Value = "  +1.09%  "
